# Update existing rows 2-13 and add new rows 14-17 to reflect new TPM-based calculations
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Mdk"
$ws.Cells.Item(2,3).Value = "Ptprz1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = [double]"3"
$ws.Cells.Item(2,6).Value = [double]"1"
$ws.Cells.Item(2,7).Value = [double]"2.180165333333334"
$ws.Cells.Item(2,8).Value = [double]"6.540496"
$ws.Cells.Item(2,9).Value = [double]"0.01970539991828544"
$ws.Cells.Item(2,10).Value = [double]"0.01970539991828544"
$ws.Cells.Item(2,11).Value = [double]"2"
$ws.Cells.Item(2,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(2,13).Value = [double]"0.138766"
$ws.Cells.Item(2,14).Value = [double]"0.416298"
$ws.Cells.Item(2,15).Value = [double]"0.01356925767068476"
$ws.Cells.Item(2,16).Value = [double]"0.01356925767068476"
$ws.Cells.Item(2,17).Value = [double]"0.3025328226453334"
$ws.Cells.Item(2,18).Value = [double]"2.722795403808"
$ws.Cells.Item(2,19).Value = [double]"0.0002673876489951054"
$ws.Cells.Item(2,20).Value = [double]"0.0002673876489951055"

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Mdk"
$ws.Cells.Item(3,3).Value = "Ptprz1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = [double]"3"
$ws.Cells.Item(3,6).Value = [double]"1"
$ws.Cells.Item(3,7).Value = [double]"2.180165333333334"
$ws.Cells.Item(3,8).Value = [double]"6.540496"
$ws.Cells.Item(3,9).Value = [double]"0.01970539991828544"
$ws.Cells.Item(3,10).Value = [double]"0.01970539991828544"
$ws.Cells.Item(3,11).Value = [double]"2"
$ws.Cells.Item(3,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(3,13).Value = [double]"0.05416133333333333"
$ws.Cells.Item(3,14).Value = [double]"0.162484"
$ws.Cells.Item(3,15).Value = [double]"0.00529617548814441"
$ws.Cells.Item(3,16).Value = [double]"0.005296175488144411"
$ws.Cells.Item(3,17).Value = [double]"0.1180806613404445"
$ws.Cells.Item(3,18).Value = [double]"1.062725952064"
$ws.Cells.Item(3,19).Value = [double]"0.0001043632560313062"
$ws.Cells.Item(3,20).Value = [double]"0.0001043632560313062"

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Mdk"
$ws.Cells.Item(4,3).Value = "Ptprz1"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = [double]"3"
$ws.Cells.Item(4,6).Value = [double]"1"
$ws.Cells.Item(4,7).Value = [double]"2.180165333333334"
$ws.Cells.Item(4,8).Value = [double]"6.540496"
$ws.Cells.Item(4,9).Value = [double]"0.01970539991828544"
$ws.Cells.Item(4,10).Value = [double]"0.01970539991828544"
$ws.Cells.Item(4,11).Value = [double]"3"
$ws.Cells.Item(4,12).Value = [double]"1"
$ws.Cells.Item(4,13).Value = [double]"10.002366"
$ws.Cells.Item(4,14).Value = [double]"30.007098"
$ws.Cells.Item(4,15).Value = [double]"0.9780831152479456"
$ws.Cells.Item(4,16).Value = [double]"0.9780831152479456"
$ws.Cells.Item(4,17).Value = [double]"21.806811604512"
$ws.Cells.Item(4,18).Value = [double]"196.261304440608"
$ws.Cells.Item(4,19).Value = [double]"0.01927351893928323"
$ws.Cells.Item(4,20).Value = [double]"0.01927351893928323"

# Row 5
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Mdk"
$ws.Cells.Item(5,3).Value = "Ptprz1"
$ws.Cells.Item(5,4).Value = "Resolving-Mac"
$ws.Cells.Item(5,5).Value = [double]"3"
$ws.Cells.Item(5,6).Value = [double]"1"
$ws.Cells.Item(5,7).Value = [double]"2.180165333333334"
$ws.Cells.Item(5,8).Value = [double]"6.540496"
$ws.Cells.Item(5,9).Value = [double]"0.01970539991828544"
$ws.Cells.Item(5,10).Value = [double]"0.01970539991828544"
$ws.Cells.Item(5,11).Value = [double]"1"
$ws.Cells.Item(5,12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(5,13).Value = [double]"0.03120566666666667"
$ws.Cells.Item(5,14).Value = [double]"0.093617"
$ws.Cells.Item(5,15).Value = [double]"0.003051451593225274"
$ws.Cells.Item(5,16).Value = [double]"0.003051451593225274"
$ws.Cells.Item(5,17).Value = [double]"0.06803351267022223"
$ws.Cells.Item(5,18).Value = [double]"0.612301614032"
$ws.Cells.Item(5,19).Value = [double]"6.013007397579327e-05"
$ws.Cells.Item(5,20).Value = [double]"6.013007397579327e-05"

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Mdk"
$ws.Cells.Item(6,3).Value = "Ptprz1"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = [double]"3"
$ws.Cells.Item(6,6).Value = [double]"1"
$ws.Cells.Item(6,7).Value = [double]"81.17653533333333"
$ws.Cells.Item(6,8).Value = [double]"243.529606"
$ws.Cells.Item(6,9).Value = [double]"0.733713204346044"
$ws.Cells.Item(6,10).Value = [double]"0.7337132043460441"
$ws.Cells.Item(6,11).Value = [double]"2"
$ws.Cells.Item(6,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(6,13).Value = [double]"0.138766"
$ws.Cells.Item(6,14).Value = [double]"0.416298"
$ws.Cells.Item(6,15).Value = [double]"0.01356925767068476"
$ws.Cells.Item(6,16).Value = [double]"0.01356925767068476"
$ws.Cells.Item(6,17).Value = [double]"11.26454310206533"
$ws.Cells.Item(6,18).Value = [double]"101.380887918588"
$ws.Cells.Item(6,19).Value = [double]"0.009955943526155251"
$ws.Cells.Item(6,20).Value = [double]"0.009955943526155253"

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Mdk"
$ws.Cells.Item(7,3).Value = "Ptprz1"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = [double]"3"
$ws.Cells.Item(7,6).Value = [double]"1"
$ws.Cells.Item(7,7).Value = [double]"81.17653533333333"
$ws.Cells.Item(7,8).Value = [double]"243.529606"
$ws.Cells.Item(7,9).Value = [double]"0.733713204346044"
$ws.Cells.Item(7,10).Value = [double]"0.7337132043460441"
$ws.Cells.Item(7,11).Value = [double]"2"
$ws.Cells.Item(7,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(7,13).Value = [double]"0.05416133333333333"
$ws.Cells.Item(7,14).Value = [double]"0.162484"
$ws.Cells.Item(7,15).Value = [double]"0.00529617548814441"
$ws.Cells.Item(7,16).Value = [double]"0.005296175488144411"
$ws.Cells.Item(7,17).Value = [double]"4.396629389033778"
$ws.Cells.Item(7,18).Value = [double]"39.569664501304"
$ws.Cells.Item(7,19).Value = [double]"0.00388587388818541"
$ws.Cells.Item(7,20).Value = [double]"0.00388587388818541"

# Row 8
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Mdk"
$ws.Cells.Item(8,3).Value = "Ptprz1"
$ws.Cells.Item(8,4).Value = "MuSCs"
$ws.Cells.Item(8,5).Value = [double]"3"
$ws.Cells.Item(8,6).Value = [double]"1"
$ws.Cells.Item(8,7).Value = [double]"81.17653533333333"
$ws.Cells.Item(8,8).Value = [double]"243.529606"
$ws.Cells.Item(8,9).Value = [double]"0.733713204346044"
$ws.Cells.Item(8,10).Value = [double]"0.7337132043460441"
$ws.Cells.Item(8,11).Value = [double]"3"
$ws.Cells.Item(8,12).Value = [double]"1"
$ws.Cells.Item(8,13).Value = [double]"10.002366"
$ws.Cells.Item(8,14).Value = [double]"30.007098"
$ws.Cells.Item(8,15).Value = [double]"0.9780831152479456"
$ws.Cells.Item(8,16).Value = [double]"0.9780831152479456"
$ws.Cells.Item(8,17).Value = [double]"811.957417015932"
$ws.Cells.Item(8,18).Value = [double]"7307.616753143388"
$ws.Cells.Item(8,19).Value = [double]"0.7176324966053312"
$ws.Cells.Item(8,20).Value = [double]"0.7176324966053313"

# Row 9
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Mdk"
$ws.Cells.Item(9,3).Value = "Ptprz1"
$ws.Cells.Item(9,4).Value = "Resolving-Mac"
$ws.Cells.Item(9,5).Value = [double]"3"
$ws.Cells.Item(9,6).Value = [double]"1"
$ws.Cells.Item(9,7).Value = [double]"81.17653533333333"
$ws.Cells.Item(9,8).Value = [double]"243.529606"
$ws.Cells.Item(9,9).Value = [double]"0.733713204346044"
$ws.Cells.Item(9,10).Value = [double]"0.7337132043460441"
$ws.Cells.Item(9,11).Value = [double]"1"
$ws.Cells.Item(9,12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(9,13).Value = [double]"0.03120566666666667"
$ws.Cells.Item(9,14).Value = [double]"0.093617"
$ws.Cells.Item(9,15).Value = [double]"0.003051451593225274"
$ws.Cells.Item(9,16).Value = [double]"0.003051451593225274"
$ws.Cells.Item(9,17).Value = [double]"2.533167902766889"
$ws.Cells.Item(9,18).Value = [double]"22.798511124902"
$ws.Cells.Item(9,19).Value = [double]"0.002238890326372157"
$ws.Cells.Item(9,20).Value = [double]"0.002238890326372157"

# Row 10
$ws.Cells.Item(10,1).Value = "MuSCs"
$ws.Cells.Item(10,2).Value = "Mdk"
$ws.Cells.Item(10,3).Value = "Ptprz1"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = [double]"3"
$ws.Cells.Item(10,6).Value = [double]"1"
$ws.Cells.Item(10,7).Value = [double]"25.672264"
$ws.Cells.Item(10,8).Value = [double]"77.016792"
$ws.Cells.Item(10,9).Value = [double]"0.2320384702908474"
$ws.Cells.Item(10,10).Value = [double]"0.2320384702908474"
$ws.Cells.Item(10,11).Value = [double]"2"
$ws.Cells.Item(10,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(10,13).Value = [double]"0.138766"
$ws.Cells.Item(10,14).Value = [double]"0.416298"
$ws.Cells.Item(10,15).Value = [double]"0.01356925767068476"
$ws.Cells.Item(10,16).Value = [double]"0.01356925767068476"
$ws.Cells.Item(10,17).Value = [double]"3.562437386224"
$ws.Cells.Item(10,18).Value = [double]"32.061936476016"
$ws.Cells.Item(10,19).Value = [double]"0.003148589792888038"
$ws.Cells.Item(10,20).Value = [double]"0.003148589792888038"

# Row 11
$ws.Cells.Item(11,1).Value = "MuSCs"
$ws.Cells.Item(11,2).Value = "Mdk"
$ws.Cells.Item(11,3).Value = "Ptprz1"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = [double]"3"
$ws.Cells.Item(11,6).Value = [double]"1"
$ws.Cells.Item(11,7).Value = [double]"25.672264"
$ws.Cells.Item(11,8).Value = [double]"77.016792"
$ws.Cells.Item(11,9).Value = [double]"0.2320384702908474"
$ws.Cells.Item(11,10).Value = [double]"0.2320384702908474"
$ws.Cells.Item(11,11).Value = [double]"2"
$ws.Cells.Item(11,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(11,13).Value = [double]"0.05416133333333333"
$ws.Cells.Item(11,14).Value = [double]"0.162484"
$ws.Cells.Item(11,15).Value = [double]"0.00529617548814441"
$ws.Cells.Item(11,16).Value = [double]"0.005296175488144411"
$ws.Cells.Item(11,17).Value = [double]"1.390444047925333"
$ws.Cells.Item(11,18).Value = [double]"12.513996431328"
$ws.Cells.Item(11,19).Value = [double]"0.001228916458660911"
$ws.Cells.Item(11,20).Value = [double]"0.001228916458660911"

# Row 12
$ws.Cells.Item(12,1).Value = "MuSCs"
$ws.Cells.Item(12,2).Value = "Mdk"
$ws.Cells.Item(12,3).Value = "Ptprz1"
$ws.Cells.Item(12,4).Value = "MuSCs"
$ws.Cells.Item(12,5).Value = [double]"3"
$ws.Cells.Item(12,6).Value = [double]"1"
$ws.Cells.Item(12,7).Value = [double]"25.672264"
$ws.Cells.Item(12,8).Value = [double]"77.016792"
$ws.Cells.Item(12,9).Value = [double]"0.2320384702908474"
$ws.Cells.Item(12,10).Value = [double]"0.2320384702908474"
$ws.Cells.Item(12,11).Value = [double]"3"
$ws.Cells.Item(12,12).Value = [double]"1"
$ws.Cells.Item(12,13).Value = [double]"10.002366"
$ws.Cells.Item(12,14).Value = [double]"30.007098"
$ws.Cells.Item(12,15).Value = [double]"0.9780831152479456"
$ws.Cells.Item(12,16).Value = [double]"0.9780831152479456"
$ws.Cells.Item(12,17).Value = [double]"256.783380576624"
$ws.Cells.Item(12,18).Value = [double]"2311.050425189616"
$ws.Cells.Item(12,19).Value = [double]"0.2269529098794399"
$ws.Cells.Item(12,20).Value = [double]"0.2269529098794399"

# Row 13
$ws.Cells.Item(13,1).Value = "MuSCs"
$ws.Cells.Item(13,2).Value = "Mdk"
$ws.Cells.Item(13,3).Value = "Ptprz1"
$ws.Cells.Item(13,4).Value = "Resolving-Mac"
$ws.Cells.Item(13,5).Value = [double]"3"
$ws.Cells.Item(13,6).Value = [double]"1"
$ws.Cells.Item(13,7).Value = [double]"25.672264"
$ws.Cells.Item(13,8).Value = [double]"77.016792"
$ws.Cells.Item(13,9).Value = [double]"0.2320384702908474"
$ws.Cells.Item(13,10).Value = [double]"0.2320384702908474"
$ws.Cells.Item(13,11).Value = [double]"1"
$ws.Cells.Item(13,12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(13,13).Value = [double]"0.03120566666666667"
$ws.Cells.Item(13,14).Value = [double]"0.093617"
$ws.Cells.Item(13,15).Value = [double]"0.003051451593225274"
$ws.Cells.Item(13,16).Value = [double]"0.003051451593225274"
$ws.Cells.Item(13,17).Value = [double]"0.8011201129626667"
$ws.Cells.Item(13,18).Value = [double]"7.210081016664"
$ws.Cells.Item(13,19).Value = [double]"0.0007080541598585617"
$ws.Cells.Item(13,20).Value = [double]"0.0007080541598585617"

# Row 14
$ws.Cells.Item(14,1).Value = "Resolving-Mac"
$ws.Cells.Item(14,2).Value = "Mdk"
$ws.Cells.Item(14,3).Value = "Ptprz1"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = [double]"3"
$ws.Cells.Item(14,6).Value = [double]"1"
$ws.Cells.Item(14,7).Value = [double]"1.608999666666667"
$ws.Cells.Item(14,8).Value = [double]"4.826999"
$ws.Cells.Item(14,9).Value = [double]"0.01454292544482312"
$ws.Cells.Item(14,10).Value = [double]"0.01454292544482312"
$ws.Cells.Item(14,11).Value = [double]"2"
$ws.Cells.Item(14,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(14,13).Value = [double]"0.138766"
$ws.Cells.Item(14,14).Value = [double]"0.416298"
$ws.Cells.Item(14,15).Value = [double]"0.01356925767068476"
$ws.Cells.Item(14,16).Value = [double]"0.01356925767068476"
$ws.Cells.Item(14,17).Value = [double]"0.2232744477446666"
$ws.Cells.Item(14,18).Value = [double]"2.009470029702"
$ws.Cells.Item(14,19).Value = [double]"0.0001973367026463627"
$ws.Cells.Item(14,20).Value = [double]"0.0001973367026463627"

# Row 15
$ws.Cells.Item(15,1).Value = "Resolving-Mac"
$ws.Cells.Item(15,2).Value = "Mdk"
$ws.Cells.Item(15,3).Value = "Ptprz1"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = [double]"3"
$ws.Cells.Item(15,6).Value = [double]"1"
$ws.Cells.Item(15,7).Value = [double]"1.608999666666667"
$ws.Cells.Item(15,8).Value = [double]"4.826999"
$ws.Cells.Item(15,9).Value = [double]"0.01454292544482312"
$ws.Cells.Item(15,10).Value = [double]"0.01454292544482312"
$ws.Cells.Item(15,11).Value = [double]"2"
$ws.Cells.Item(15,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(15,13).Value = [double]"0.05416133333333333"
$ws.Cells.Item(15,14).Value = [double]"0.162484"
$ws.Cells.Item(15,15).Value = [double]"0.00529617548814441"
$ws.Cells.Item(15,16).Value = [double]"0.005296175488144411"
$ws.Cells.Item(15,17).Value = [double]"0.08714556727955555"
$ws.Cells.Item(15,18).Value = [double]"0.7843101055159999"
$ws.Cells.Item(15,19).Value = [double]"7.702188526678388e-05"
$ws.Cells.Item(15,20).Value = [double]"7.702188526678389e-05"

# Row 16
$ws.Cells.Item(16,1).Value = "Resolving-Mac"
$ws.Cells.Item(16,2).Value = "Mdk"
$ws.Cells.Item(16,3).Value = "Ptprz1"
$ws.Cells.Item(16,4).Value = "MuSCs"
$ws.Cells.Item(16,5).Value = [double]"3"
$ws.Cells.Item(16,6).Value = [double]"1"
$ws.Cells.Item(16,7).Value = [double]"1.608999666666667"
$ws.Cells.Item(16,8).Value = [double]"4.826999"
$ws.Cells.Item(16,9).Value = [double]"0.01454292544482312"
$ws.Cells.Item(16,10).Value = [double]"0.01454292544482312"
$ws.Cells.Item(16,11).Value = [double]"3"
$ws.Cells.Item(16,12).Value = [double]"1"
$ws.Cells.Item(16,13).Value = [double]"10.002366"
$ws.Cells.Item(16,14).Value = [double]"30.007098"
$ws.Cells.Item(16,15).Value = [double]"0.9780831152479456"
$ws.Cells.Item(16,16).Value = [double]"0.9780831152479456"
$ws.Cells.Item(16,17).Value = [double]"16.093803559878"
$ws.Cells.Item(16,18).Value = [double]"144.844232038902"
$ws.Cells.Item(16,19).Value = [double]"0.01422418982389122"
$ws.Cells.Item(16,20).Value = [double]"0.01422418982389122"

# Row 17
$ws.Cells.Item(17,1).Value = "Resolving-Mac"
$ws.Cells.Item(17,2).Value = "Mdk"
$ws.Cells.Item(17,3).Value = "Ptprz1"
$ws.Cells.Item(17,4).Value = "Resolving-Mac"
$ws.Cells.Item(17,5).Value = [double]"3"
$ws.Cells.Item(17,6).Value = [double]"1"
$ws.Cells.Item(17,7).Value = [double]"1.608999666666667"
$ws.Cells.Item(17,8).Value = [double]"4.826999"
$ws.Cells.Item(17,9).Value = [double]"0.01454292544482312"
$ws.Cells.Item(17,10).Value = [double]"0.01454292544482312"
$ws.Cells.Item(17,11).Value = [double]"1"
$ws.Cells.Item(17,12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(17,13).Value = [double]"0.03120566666666667"
$ws.Cells.Item(17,14).Value = [double]"0.093617"
$ws.Cells.Item(17,15).Value = [double]"0.003051451593225274"
$ws.Cells.Item(17,16).Value = [double]"0.003051451593225274"
$ws.Cells.Item(17,17).Value = [double]"0.05020990726477778"
$ws.Cells.Item(17,18).Value = [double]"0.451889165383"
$ws.Cells.Item(17,19).Value = [double]"4.43770330187619e-05"
$ws.Cells.Item(17,20).Value = [double]"4.43770330187619e-05"

Write-Host "Updated sheet to new TPM values; dimension now A1:T17"